$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1994.132995605469
$ws.Range("C3").Value = 3989.934921264648
$ws.Range("D3").Value = 6980.89599609375
$ws.Range("E3").Value = 9972.572326660156
$ws.Range("F3").Value = 12964.24865722656
$ws.Range("G3").Value = 15957.35549926758
$ws.Range("H3").Value = 19946.81358337402
$ws.Range("I3").Value = 22938.48991394043
$ws.Range("J3").Value = 25896.78764343262
$ws.Range("K3").Value = 30897.37892150879
$ws.Range("B4").Value = 2992.153167724609
$ws.Range("C4").Value = 4987.955093383789
$ws.Range("D4").Value = 6981.611251831055
$ws.Range("E4").Value = 8976.459503173828
$ws.Range("F4").Value = 11968.37425231934
$ws.Range("G4").Value = 15959.02442932129
$ws.Range("H4").Value = 16954.42199707031
$ws.Range("I4").Value = 19947.05200195312
$ws.Range("J4").Value = 23962.97454833984
$ws.Range("K4").Value = 25928.97415161133
$ws.Range("B5").Value = 2991.914749145508
$ws.Range("C5").Value = 6980.419158935547
$ws.Range("D5").Value = 11969.32792663574
$ws.Range("E5").Value = 16954.18357849121
$ws.Range("F5").Value = 20943.88008117676
$ws.Range("G5").Value = 25929.92782592773
$ws.Range("H5").Value = 30917.16766357422
$ws.Range("I5").Value = 35904.4075012207
$ws.Range("J5").Value = 39864.77851867676
$ws.Range("K5").Value = 45877.45666503906
$ws.Range("B6").Value = 27925.01449584961
$ws.Range("C6").Value = 115660.4290008545
$ws.Range("D6").Value = 250365.7341003418
$ws.Range("E6").Value = 449797.3918914795
$ws.Range("F6").Value = 701125.8602142334
$ws.Range("G6").Value = 1031243.085861206
$ws.Range("H6").Value = 1378316.402435303
$ws.Range("I6").Value = 1795199.632644653
$ws.Range("J6").Value = 2350753.545761108
$ws.Range("K6").Value = 2811478.137969971
$ws.Range("B7").Value = 43919.80171203613
$ws.Range("C7").Value = 151596.5461730957
$ws.Range("D7").Value = 360000.8487701416
$ws.Range("E7").Value = 666218.2807922363
$ws.Range("F7").Value = 998328.447341919
$ws.Range("G7").Value = 1890941.858291626
$ws.Range("H7").Value = 2160223.245620728
$ws.Range("I7").Value = 2740709.066390991
$ws.Range("J7").Value = 3216392.278671265
$ws.Range("K7").Value = 3729995.012283325
